$d = $word.ActiveDocument

function Get-ParaIndexByText {
    param($doc, [string]$targetText)
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $i
        }
    }
    return -1
}

# --- Change 1 -------------------------------------------------------------
# After "Crear la BBDD y conectar con el BackEnd." add a new bullet:
# "Crear las clases en el paquete entidades (BackEnd)"
$idx1 = Get-ParaIndexByText $d "Crear la BBDD y conectar con el BackEnd."
$p1 = $d.Paragraphs.Item($idx1)
$p1.Range.InsertParagraphAfter()

$idx1b = Get-ParaIndexByText $d "Crear la BBDD y conectar con el BackEnd."
$newP1 = $d.Paragraphs.Item($idx1b + 1)
$newP1.Range.Text = "Crear las clases en el paquete entidades (BackEnd)"

# --- Change 2 -------------------------------------------------------------
# Rewrite "Crear los paquetes y clases necesarios para el proyecto." as
# "Crear los Servicios y las Interfaces para cada clase. (BackEnd)"
$idx2 = Get-ParaIndexByText $d "Crear los paquetes y clases necesarios para el proyecto."
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range
[void]$r2.MoveEnd(1, -1)
$r2.Text = "Crear los Servicios y las Interfaces para cada clase. (BackEnd)"

# Then add a new bullet right after it:
# "Crear los Controlador para cada clase. (BackEnd)"
$idx2b = Get-ParaIndexByText $d "Crear los Servicios y las Interfaces para cada clase. (BackEnd)"
$p2b = $d.Paragraphs.Item($idx2b)
$p2b.Range.InsertParagraphAfter()

$idx2c = Get-ParaIndexByText $d "Crear los Servicios y las Interfaces para cada clase. (BackEnd)"
$newP2 = $d.Paragraphs.Item($idx2c + 1)
$newP2.Range.Text = "Crear los Controlador para cada clase. (BackEnd)"

Write-Output "done"
